# Insert a new data row at row 80 (pushing existing rows 80:145 down to 81:146,
# and growing the sheet from A1:T145 to A1:T146), then populate the new row
# with the new Papaya price-report entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 80 (and everything below it) down by one row.
$ws.Rows.Item(80).EntireRow.Insert()

# Fill in the newly inserted row 80 with the new record's data.
$ws.Cells.Item(80, 1).Value = 10
$ws.Cells.Item(80, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(80, 3).Value = "La Araucanía"
$ws.Cells.Item(80, 4).Value = 45264
$ws.Cells.Item(80, 5).Value = 9
$ws.Cells.Item(80, 6).Value = "Fruta"
$ws.Cells.Item(80, 7).Value = 100108
$ws.Cells.Item(80, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(80, 9).Value = 100108004
$ws.Cells.Item(80, 10).Value = "Papaya"
$ws.Cells.Item(80, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(80, 12).Value = "Primera"
$ws.Cells.Item(80, 13).Value = 40
$ws.Cells.Item(80, 14).Value = 37000
$ws.Cells.Item(80, 15).Value = 37000
$ws.Cells.Item(80, 16).Value = 37000
$ws.Cells.Item(80, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(80, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(80, 19).Value = 2467
$ws.Cells.Item(80, 20).Value = 15
